# Add the new day's row of data (row 19) to the 100 Error Counts log
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = 45961
$ws.Range("B19").Value = 594
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 579

# Move/leave the active selection where the user last clicked after entering data
$ws.Range("E27").Select() | Out-Null
